$d = $word.ActiveDocument
$d.Content.Find.Execute("solidfy", $true, $false, $false, $false, $false,
                         $true, 1, $false, "solidify", 2)
